# Auto-generated edit script applying the diff to Sheets/Sagittarius_Profits.xlsx
# Updates currentAveragePrice / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ
# values (columns H-N) for specific leve rows across sheets, per the scheduled runner diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 31
$ws.Range("H31").Value = 1861.25
$ws.Range("I31").Value = 1861.25
$ws.Range("K31").Value = 5583.75
$ws.Range("M31").Value = -5353.75
# Row 32
$ws.Range("H32").Value = 3408.2307
$ws.Range("I32").Value = 878.6667
$ws.Range("J32").Value = 4167.1
$ws.Range("K32").Value = 878.6667
$ws.Range("L32").Value = 4167.1
$ws.Range("M32").Value = -552.6667
$ws.Range("N32").Value = -4819.1
# Row 39
$ws.Range("H39").Value = 1333.5
$ws.Range("I39").Value = 522.7143
$ws.Range("K39").Value = 1568.1429
$ws.Range("M39").Value = -1272.1429
# Row 58
$ws.Range("H58").Value = 430.81818
$ws.Range("J58").Value = 505.7143
$ws.Range("L58").Value = 1517.1429
$ws.Range("N58").Value = -1817.1429

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1207.9
$ws.Range("I45").Value = 1219.8889
$ws.Range("J45").Value = 1100
$ws.Range("K45").Value = 1219.8889
$ws.Range("L45").Value = 1100
$ws.Range("M45").Value = -842.8888999999999
$ws.Range("N45").Value = -1854
# Row 60
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
# Row 61
$ws.Range("H61").Value = 1949
$ws.Range("I61").Value = 1937.7778
$ws.Range("K61").Value = 1937.7778
$ws.Range("M61").Value = -1725.7778
# Row 74
$ws.Range("H74").Value = 3070.8572
$ws.Range("I74").Value = 3339.8
$ws.Range("J74").Value = 2398.5
$ws.Range("K74").Value = 3339.8
$ws.Range("L74").Value = 2398.5
$ws.Range("M74").Value = -2465.8
$ws.Range("N74").Value = -4146.5
# Row 77
$ws.Range("H77").Value = 3070.8572
$ws.Range("I77").Value = 3339.8
$ws.Range("J77").Value = 2398.5
$ws.Range("K77").Value = 16699
$ws.Range("L77").Value = 11992.5
$ws.Range("M77").Value = -12331
$ws.Range("N77").Value = -20728.5
# Row 112
$ws.Range("H112").Value = 28333.334
$ws.Range("J112").Value = 28333.334
$ws.Range("L112").Value = 28333.334
$ws.Range("N112").Value = -31287.334
# Row 136
$ws.Range("H136").Value = 1949
$ws.Range("I136").Value = 1937.7778
$ws.Range("K136").Value = 5813.3334
$ws.Range("M136").Value = -3263.3334

$ws = $wb.Worksheets.Item("BSM")
# Row 75
$ws.Range("H75").Value = 33177.6
$ws.Range("I75").Value = 33177.6
$ws.Range("K75").Value = 33177.6
$ws.Range("M75").Value = -32241.6
# Row 76
$ws.Range("H76").Value = 76939.14
$ws.Range("J76").Value = 76939.14
$ws.Range("L76").Value = 76939.14
$ws.Range("N76").Value = -77569.14
# Row 78
$ws.Range("H78").Value = 33177.6
$ws.Range("I78").Value = 33177.6
$ws.Range("K78").Value = 99532.79999999999
$ws.Range("M78").Value = -94852.79999999999
# Row 79
$ws.Range("H79").Value = 76939.14
$ws.Range("J79").Value = 76939.14
$ws.Range("L79").Value = 76939.14
$ws.Range("N79").Value = -79123.14
# Row 92
$ws.Range("H92").Value = 38749.332
$ws.Range("J92").Value = 38749.332
$ws.Range("L92").Value = 38749.332
$ws.Range("N92").Value = -43741.332
# Row 97
$ws.Range("H97").Value = 18609.5
$ws.Range("I97").Value = 18179.334
$ws.Range("K97").Value = 18179.334
$ws.Range("M97").Value = -17188.334
# Row 107
$ws.Range("H107").Value = 401.08334
$ws.Range("I107").Value = 383.54544
$ws.Range("K107").Value = 383.54544
$ws.Range("M107").Value = 1536.45456
# Row 134
$ws.Range("H134").Value = 2040.3636
$ws.Range("I134").Value = 2040.3636
$ws.Range("K134").Value = 6121.0908
$ws.Range("M134").Value = -3586.0908

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 53.333332
$ws.Range("I7").Value = 55
$ws.Range("J7").Value = 50
$ws.Range("K7").Value = 55
$ws.Range("L7").Value = 50
$ws.Range("M7").Value = 58
$ws.Range("N7").Value = -276
# Row 58
$ws.Range("H58").Value = 2117.4285
$ws.Range("I58").Value = 1914.8334
$ws.Range("K58").Value = 1914.8334
$ws.Range("M58").Value = -1711.8334
# Row 107
$ws.Range("H107").Value = 646.6
$ws.Range("I107").Value = 316.5
$ws.Range("K107").Value = 316.5
$ws.Range("M107").Value = 1603.5
# Row 132
$ws.Range("H132").Value = 2733.7856
$ws.Range("I132").Value = 2681.4167
$ws.Range("J132").Value = 3048
$ws.Range("K132").Value = 8044.250100000001
$ws.Range("L132").Value = 9144
$ws.Range("M132").Value = -5514.250100000001
$ws.Range("N132").Value = -14204
# Row 134
$ws.Range("H134").Value = 2390.9412
$ws.Range("I134").Value = 2401.5
$ws.Range("K134").Value = 7204.5
$ws.Range("M134").Value = -4669.5
# Row 136
$ws.Range("H136").Value = 2117.4285
$ws.Range("I136").Value = 1914.8334
$ws.Range("K136").Value = 5744.5002
$ws.Range("M136").Value = -3194.5002

$ws = $wb.Worksheets.Item("CUL")
# Row 69
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 936.25
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 2808.75
$ws.Range("N69").Value = -4430.75
$ws.Range("M69").ClearContents()
# Row 72
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 936.25
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 8426.25
$ws.Range("N72").Value = -16538.25
$ws.Range("M72").ClearContents()
# Row 98
$ws.Range("H98").Value = 2824.6667
$ws.Range("J98").Value = 2824.6667
$ws.Range("L98").Value = 8474.000100000001
$ws.Range("N98").Value = -11470.0001
# Row 109
$ws.Range("H109").Value = 1672.5
$ws.Range("I109").Value = 846.6667
$ws.Range("J109").Value = 4150
$ws.Range("K109").Value = 2540.0001
$ws.Range("L109").Value = 12450
$ws.Range("M109").Value = -1500.0001
$ws.Range("N109").Value = -14530

$ws = $wb.Worksheets.Item("GSM")
# Row 62
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
# Row 65
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
# Row 70
$ws.Range("H70").Value = 5504
$ws.Range("I70").Value = 5008
$ws.Range("J70").Value = 6000
$ws.Range("K70").Value = 5008
$ws.Range("L70").Value = 6000
$ws.Range("M70").Value = -4738
$ws.Range("N70").Value = -6540
# Row 73
$ws.Range("H73").Value = 5504
$ws.Range("I73").Value = 5008
$ws.Range("J73").Value = 6000
$ws.Range("K73").Value = 5008
$ws.Range("L73").Value = 6000
$ws.Range("M73").Value = -4072
$ws.Range("N73").Value = -7872
# Row 111
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 1190.6
$ws.Range("I61").Value = 866
$ws.Range("J61").Value = 1677.5
$ws.Range("K61").Value = 866
$ws.Range("L61").Value = 1677.5
$ws.Range("M61").Value = -664
$ws.Range("N61").Value = -2081.5
# Row 63
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
# Row 66
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
# Row 110
$ws.Range("H110").Value = 35000
$ws.Range("J110").Value = 35000
$ws.Range("L110").Value = 35000
$ws.Range("N110").Value = -43180
# Row 113
$ws.Range("H113").Value = 1190.6
$ws.Range("I113").Value = 866
$ws.Range("J113").Value = 1677.5
$ws.Range("K113").Value = 866
$ws.Range("L113").Value = 1677.5
$ws.Range("M113").Value = 1304
$ws.Range("N113").Value = -6017.5

$ws = $wb.Worksheets.Item("WVR")
# Row 70
$ws.Range("H70").Value = 67287
$ws.Range("J70").Value = 67287
$ws.Range("L70").Value = 67287
$ws.Range("N70").Value = -67917
# Row 73
$ws.Range("H73").Value = 67287
$ws.Range("J73").Value = 67287
$ws.Range("L73").Value = 67287
$ws.Range("N73").Value = -69471
# Row 75
$ws.Range("H75").Value = 76501
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()
# Row 78
$ws.Range("H78").Value = 76501
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()
# Row 107
$ws.Range("H107").Value = 396.66666
$ws.Range("I107").Value = 200.66667
$ws.Range("J107").Value = 592.6667
$ws.Range("K107").Value = 602.00001
$ws.Range("L107").Value = 1778.0001
$ws.Range("M107").Value = 1317.99999
$ws.Range("N107").Value = -5618.0001
# Row 136
$ws.Range("H136").Value = 2293.2942
$ws.Range("J136").Value = 1366.3334
$ws.Range("L136").Value = 4099.0002
$ws.Range("N136").Value = -9199.0002

